# Atualização de bases das ligas, do dia: 29-02-2024 às 07:50
# Rows 83 and 84 (game id 81 and 82) had their match data swapped between
# each other - everything except the row index in column A stays attached
# to the same physical row, but "id" (B) plus all the match/odds columns
# (F:AC) now belong to the other row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 83 (id=81) now carries what used to be row 84's data ----
$ws.Range("B83").Value = 6227884
$ws.Range("F83").Value = "Cavalry FC"
$ws.Range("G83").Value = "Pacific FC CA"
$ws.Range("H83").Value = 3
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = "H"
$ws.Range("K83").Value = 2.25
$ws.Range("L83").Value = 3.1
$ws.Range("M83").Value = 2.875
$ws.Range("N83").Value = 2.05
$ws.Range("O83").Value = 3.2
$ws.Range("P83").Value = 3.2
$ws.Range("Q83").Value = -0.25
$ws.Range("R83").Value = 1.825
$ws.Range("S83").Value = 1.975
$ws.Range("T83").Value = 2.5
$ws.Range("U83").Value = 1.825
$ws.Range("V83").Value = 1.975
$ws.Range("W83").Value = 1.05
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = 0.825
$ws.Range("AA83").Value = -1
$ws.Range("AB83").Value = 0.825
$ws.Range("AC83").Value = -1

# ---- Row 84 (id=82) now carries what used to be row 83's data ----
$ws.Range("B84").Value = 7301364
$ws.Range("F84").Value = "Forge FC"
$ws.Range("G84").Value = "Atletico Ottawa"
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = "A"
$ws.Range("K84").Value = 1.8
$ws.Range("L84").Value = 3.6
$ws.Range("M84").Value = 3.5
$ws.Range("N84").Value = 1.533
$ws.Range("O84").Value = 3.8
$ws.Range("P84").Value = 5
$ws.Range("Q84").Value = -1
$ws.Range("R84").Value = 1.975
$ws.Range("S84").Value = 1.825
$ws.Range("T84").Value = 2.5
$ws.Range("U84").Value = 1.9
$ws.Range("V84").Value = 1.9
$ws.Range("W84").Value = -1
$ws.Range("X84").Value = -1
$ws.Range("Y84").Value = 4
$ws.Range("Z84").Value = -1
$ws.Range("AA84").Value = 0.825
$ws.Range("AB84").Value = -1
$ws.Range("AC84").Value = 0.8999999999999999
